# Append new case-data rows (895-905) to Sheet1, as described in the commit:
# "Warnings for all add conditions on Jail and No Jail in place."
#
# Columns: A Case#, B Judge, C Charge, D Statute, E Class, F Plea, G Finding,
#          H Fine, I Fine2, J Jail, K No Jail

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row=895; A="22CRB00136"; B="Hemmeter"; C="DOMESTIC VIOLENCE";                     D="2919.25(A)";     E="No Data"; F="No Contest"; G="Guilty"; H="$ 50"; I="$ 25"; J="10";   K="None" },
    @{ Row=896; A="22CRB00136"; B="Hemmeter"; C="ASSAULT - M1";                          D="2903.13(A)";     E="No Data"; F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=897; A="21CRB01291"; B="Hemmeter"; C="PERMISSION REQ'D TO USE LICENSED DOCK"; D="1501:46-12-04";  E="MM";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J=$null;  K=$null },
    @{ Row=898; A="21CRB01291"; B="Hemmeter"; C="PERMISSION REQ'D TO USE LICENSED DOCK"; D="1501:46-12-04";  E="MM";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=899; A="21CRB01291"; B="Hemmeter"; C="PERMISSION REQ'D TO USE LICENSED DOCK"; D="1501:46-12-04";  E="MM";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=900; A="21CRB01268"; B="Hemmeter"; C="POSSESSION DRUG PARAPHERNALIA";         D="2925.14(C)";     E="M4";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=901; A="21TRD09386"; B="Hemmeter"; C="DUS UCM";                               D="4510.111";       E="UCM";     F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J=$null;  K=$null },
    @{ Row=902; A="21TRD09386"; B="Hemmeter"; C="TAIL LIGHTS-REAR LICENSE PLATE";        D="4513.05";        E="MM";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J=$null;  K=$null },
    @{ Row=903; A="21TRD09386"; B="Hemmeter"; C="DUS UCM";                               D="4510.111";       E="UCM";     F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=904; A="21TRD09386"; B="Hemmeter"; C="TAIL LIGHTS-REAR LICENSE PLATE";        D="4513.05";        E="MM";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" },
    @{ Row=905; A="03TRD13368"; B="Hemmeter"; C="SPEED REDUCED ZONE 3RD OR MORE";        D="4511.21C***";    E="M3";      F="No Contest"; G="Guilty"; H="$ 0";  I="$ 0";  J="None"; K="None" }
)

# Columns/values that look numeric to Excel (e.g. "4510.111", "$ 0", "10")
# would silently be converted from text to a numeric cell on assignment.
# To keep them as plain strings (matching the rest of the sheet) we set
# NumberFormat to Text on that single cell right before assigning it, then
# clear the temporary formatting again so no stray cell style is left
# behind. This is done per-cell (not per-row-range) so that untouched
# sibling cells are never materialized as empty <c> elements.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    Set-TextValue $ws.Cells.Item($rowNum, 4) $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    Set-TextValue $ws.Cells.Item($rowNum, 8) $r.H
    Set-TextValue $ws.Cells.Item($rowNum, 9) $r.I
    if ($r.J -ne $null) {
        Set-TextValue $ws.Cells.Item($rowNum, 10) $r.J
    }
    if ($r.K -ne $null) {
        $ws.Cells.Item($rowNum, 11).Value = $r.K
    }
}
